$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.573.32"
$ws.Range("E2").Value = "  +2.63%  "
$ws.Range("D3").Value = "2.199.15"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'258.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").Value = "'83.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.37%  "
$ws.Range("D7").Value = "'0.614"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.600"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.39%  "
$ws.Range("D10").Value = "'44.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.47%  "
$ws.Range("D11").Value = "'0.0916"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("E12").Value = "  +5.42%  "
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").Value = "2.528.83"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "'14.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "2.216.26"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "43.522.74"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "'69.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").Value = "'5.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").Value = "'2.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.13%  "
$ws.Range("D23").Value = "'231.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").Value = "'9.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.96%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +1.84%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "'39.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.44%  "
$ws.Range("E29").Value = "  +5.47%  "
$ws.Range("E30").Value = "  +3.10%  "
$ws.Range("D31").Value = "'174.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("D32").Value = "'20.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("E33").Value = "  +5.20%  "
$ws.Range("D34").Value = "'5.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.13%  "
$ws.Range("E36").Value = "  +4.87%  "
$ws.Range("D37").Value = "'4.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.79%  "
$ws.Range("E38").Value = "  +7.61%  "
$ws.Range("D39").Value = "'12.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.30%  "
$ws.Range("D40").Value = "'2.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.62%  "
$ws.Range("D41").Value = "'2.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("D42").Value = "'63.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.68%  "
$ws.Range("E43").Value = "  +6.13%  "
$ws.Range("E44").Value = "  +3.25%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'8.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'99.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("E48").Value = "  +5.09%  "
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").Value = "'0.438"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.83%  "
$ws.Range("D51").Value = "'1.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.47%  "
